# Update lexical diversity to count as proportion of all tokens, not just word tokens
# This updates the per-document lexical_diversity values on sheets "2019".."2024"
# (column B, sheet indices 1-6) and recomputes the describe()-style summary
# statistics on the "Summary" sheet (sheet index 7, columns C:I) to match.

$wb = $excel.ActiveWorkbook

# New per-document lexical_diversity values, keyed by worksheet index (1-6) -> row -> value.
$newValues = @{
    1 = @{
        2 = 0.413716814159292
        3 = 0.2957600827300931
        4 = 0.4083969465648855
        5 = 0.236078144969874
    }
    2 = @{
        2 = 0.6235294117647059
        3 = 0.4816753926701571
        4 = 0.4401408450704226
        5 = 0.4304093567251462
        6 = 0.6020942408376964
        7 = 0.4396423248882265
        8 = 0.3581213307240704
    }
    3 = @{
        2 = 0.4502487562189055
        3 = 0.3407275953859805
        4 = 0.4699140401146132
        5 = 0.6439393939393939
        6 = 0.5223613595706619
        7 = 0.4884318766066838
        8 = 0.6686746987951807
        9 = 0.5034013605442177
        10 = 0.5612244897959183
        11 = 0.5284810126582279
    }
    4 = @{
        2 = 0.4573304157549234
        3 = 0.4089012517385257
        4 = 0.6197183098591549
        5 = 0.4783783783783784
        6 = 0.3429027113237639
        7 = 0.3607784431137724
    }
    5 = @{
        2 = 0.5103626943005182
    }
    6 = @{
        2 = 0.4774011299435028
        3 = 0.564327485380117
        4 = 0.5150684931506849
    }
}

foreach ($sheetIndex in $newValues.Keys) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $rows = $newValues[$sheetIndex]
    foreach ($row in $rows.Keys) {
        $ws.Range("B$row").Value = $rows[$row]
    }
}

# Recompute the Summary sheet's per-year statistics (count, mean, std, min,
# 25%, 50%, 75%, max) from the updated per-document values, matching
# pandas'/Excel's sample statistics (std uses ddof=1, quantiles use linear
# interpolation).

$summaryRows = @{
    2 = 1
    3 = 2
    4 = 3
    5 = 4
    6 = 5
    7 = 6
}

$wsSummary = $wb.Worksheets.Item(7)

function Get-Percentile($sorted, [double]$p) {
    $n = $sorted.Count
    if ($n -eq 1) { return $sorted[0] }
    $rank = $p * ($n - 1)
    $lo = [Math]::Floor($rank)
    $hi = [Math]::Ceiling($rank)
    if ($lo -eq $hi) { return $sorted[$lo] }
    $frac = $rank - $lo
    return $sorted[$lo] + ($sorted[$hi] - $sorted[$lo]) * $frac
}

foreach ($summaryRow in $summaryRows.Keys) {
    $sheetIndex = $summaryRows[$summaryRow]
    $vals = @($newValues[$sheetIndex].Keys | Sort-Object | ForEach-Object { $newValues[$sheetIndex][$_] })
    $sorted = @($vals | Sort-Object)
    $n = $vals.Count

    $sum = 0.0
    foreach ($v in $vals) { $sum += $v }
    $mean = $sum / $n

    if ($n -gt 1) {
        $sq = 0.0
        foreach ($v in $vals) { $sq += ($v - $mean) * ($v - $mean) }
        $std = [Math]::Sqrt($sq / ($n - 1))
    } else {
        $std = $null
    }

    $min = $sorted[0]
    $max = $sorted[$n - 1]
    $p25 = Get-Percentile $sorted 0.25
    $p50 = Get-Percentile $sorted 0.50
    $p75 = Get-Percentile $sorted 0.75

    $wsSummary.Range("B$summaryRow").Value = $n
    $wsSummary.Range("C$summaryRow").Value = $mean
    if ($null -ne $std) {
        $wsSummary.Range("D$summaryRow").Value = $std
    }
    $wsSummary.Range("E$summaryRow").Value = $min
    $wsSummary.Range("F$summaryRow").Value = $p25
    $wsSummary.Range("G$summaryRow").Value = $p50
    $wsSummary.Range("H$summaryRow").Value = $p75
    $wsSummary.Range("I$summaryRow").Value = $max
}
